$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for crypto rows per latest scrape

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.473.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.36%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.756.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.26%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.17%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.77%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.757.73"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.12%  "

# Row 8
$ws.Range("E8").Value = "  -0.12%  "

# Row 9
$ws.Range("E9").Value = "  -0.82%  "

# Row 10
$ws.Range("E10").Value = "  -1.69%  "

# Row 11
$ws.Range("E11").Value = "  -0.23%  "

# Row 12
$ws.Range("E12").Value = "  -1.56%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000272"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.62%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.21%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.385.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.32%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.763.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.00%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.11%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.413.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.51%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.15%  "

# Row 20
$ws.Range("E20").Value = "  +0.67%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.37%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "466.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.77%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.716"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.71%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.26%  "

# Row 25
$ws.Range("E25").Value = "  -9.05%  "

# Row 26
$ws.Range("E26").Value = "  -1.88%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.65%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.46%  "

# Row 29
$ws.Range("E29").Value = "  -0.06%  "

# Row 30
$ws.Range("E30").Value = "  -2.42%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.903.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.29%  "

# Row 32
$ws.Range("E32").Value = "  -1.11%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.97%  "

# Row 34
$ws.Range("E34").Value = "  -4.00%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.27%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.715.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.48%  "

# Row 37
$ws.Range("E37").Value = "  +2.04%  "

# Row 38
$ws.Range("E38").Value = "  -1.34%  "

# Row 40
$ws.Range("E40").Value = "  -2.35%  "

# Row 41
$ws.Range("E41").Value = "  -2.92%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.00%  "

# Row 43
$ws.Range("E43").Value = "  -1.84%  "

# Row 45
$ws.Range("E45").Value = "  -1.08%  "

# Row 46
$ws.Range("E46").Value = "  -2.43%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.87%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "394.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.16%  "

# Row 49
$ws.Range("E49").Value = "  -6.81%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "138.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.05%  "

# Row 51
$ws.Range("E51").Value = "  -2.68%  "
